$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Mark the run holding the second inline picture (wp14:anchorId="62BC6316",
#    the screenshot showing the "drag files" step) as NoProofing, i.e. add
#    <w:rPr><w:noProof/></w:rPr> to its run.
# ---------------------------------------------------------------------------
$shp = $d.InlineShapes.Item(2)
$shp.Range.NoProofing = -1

# ---------------------------------------------------------------------------
# 2) Append the "Update" note at the very end of the document (before the
#    final section properties), as three new paragraphs:
#      "Update:"
#      "Sembra possa andare anche <i>ssh.studenti.math.unipd.it</i> come Tunnel"
#      (empty paragraph)
# ---------------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p><w:r><w:t>Update:</w:t></w:r></w:p>' +
       '<w:p>' +
       '<w:r><w:t xml:space="preserve">Sembra possa andare anche </w:t></w:r>' +
       '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>ssh.studenti.math.unipd.it</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve"> come Tunnel</w:t></w:r>' +
       '</w:p>' +
       '<w:p/>' +
       '</w:body></w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$endRange.InsertXML($xml)

Write-Host "Paragraphs after edit:" $d.Paragraphs.Count
